$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff changes cells A1:E1 on row 1 from 10 to 15
$ws.Range("A1:E1").Value = 15
